$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: keep B2 as text "1" (was "2") ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").Style = "Normal"

# --- Row 2: updated metric values ---
$ws.Range("D2").Value = 0.0933
$ws.Range("E2").Value = 0.0898
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.84
$ws.Range("J2").Value = 0.8387917808219177
$ws.Range("K2").Value = 14.6
$ws.Range("L2").Value = 0.8342857142857143
$ws.Range("M2").Value = 4.52
$ws.Range("N2").Value = 0.02574031890660592
$ws.Range("O2").Value = 0.3095890410958904
$ws.Range("P2").Value = 4.52
$ws.Range("Q2").Value = 0.02574031890660592
$ws.Range("R2").Value = 0.3095890410958904
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 1.7
$ws.Range("V2").Value = 0.00968109339407745
$ws.Range("W2").Value = 0.1546610169491525
$ws.Range("X2").Value = 0.03524549617030377
$ws.Range("Y2").Value = 0.1194155207788488
$ws.Range("Z2").Value = 0.1893324678134805
$ws.Range("AA2").Value = 0.1588105178446777
$ws.Range("AB2").Value = 0.03524549617030377
$ws.Range("AC2").Value = 0.1235650216743739
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = -1.7
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = -0.009775733179988498
$ws.Range("AK2").Value = -0.01715438950554995
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0

# --- Row 2: remove debt_ebitda/ebit_interest_expenses/net_debt_ebitda/ebit_net_interest_expenses (AN2:AQ2) ---
$ws.Range("AN2:AQ2").ClearContents()

# --- Row 3: updated metric values (Marlin Global Limited) ---
$ws.Range("D3").Value = 0.0933
$ws.Range("E3").Value = 0.0898
$ws.Range("I3").Value = 0.84
$ws.Range("J3").Value = 0.8387917808219177
$ws.Range("K3").Value = 14.6
$ws.Range("L3").Value = 0.8342857142857143
$ws.Range("M3").Value = 4.52
$ws.Range("N3").Value = 0.02574031890660592
$ws.Range("O3").Value = 0.3095890410958904
$ws.Range("P3").Value = 4.52
$ws.Range("Q3").Value = 0.02574031890660592
$ws.Range("R3").Value = 0.3095890410958904
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 1.7
$ws.Range("V3").Value = 0.00968109339407745
$ws.Range("W3").Value = 0.1546610169491525
$ws.Range("X3").Value = 0.03524549617030377
$ws.Range("Y3").Value = 0.1194155207788488
$ws.Range("Z3").Value = 0.1893324678134805
$ws.Range("AA3").Value = 0.1588105178446777
$ws.Range("AB3").Value = 0.03524549617030377
$ws.Range("AC3").Value = 0.1235650216743739
$ws.Range("AG3").Value = -1.7
$ws.Range("AJ3").Value = -0.009775733179988498
$ws.Range("AK3").Value = -0.01715438950554995

# --- Row 4 (Powerhouse Ventures Limited) removed entirely from the dataset ---
$ws.Rows(4).Delete()
